$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 810, shifting existing rows 810:866 down to 811:867
$ws.Rows.Item(810).Insert()

# Populate the new row 810 with data (mirrors the surrounding rows' structure)
$ws.Cells.Item(810, 1).Value = 9
$ws.Cells.Item(810, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(810, 3).Value = "Metropolitana"
$ws.Cells.Item(810, 4).Value = 45223
$ws.Cells.Item(810, 5).Value = 13
$ws.Cells.Item(810, 6).Value = 100112031
$ws.Cells.Item(810, 7).Value = "Poroto verde"
$ws.Cells.Item(810, 8).Value = "Magnum"
$ws.Cells.Item(810, 9).Value = "Primera"
$ws.Cells.Item(810, 10).Value = 70
$ws.Cells.Item(810, 11).Value = 26000
$ws.Cells.Item(810, 12).Value = 28000
$ws.Cells.Item(810, 13).Value = 27000
$ws.Cells.Item(810, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(810, 15).Value = "Perú"
$ws.Cells.Item(810, 16).Value = 1080
$ws.Cells.Item(810, 17).Value = 25
$ws.Cells.Item(810, 18).Value = "Hortaliza"
